# Insert a new weekly price-report row for "Vega Modelo de Temuco - Cebollín".
#
# The new record is inserted immediately after the existing row 377 (pushing
# rows 378..408 down to 379..409). The new row duplicates row 377's fixed
# attributes (market, region, category, etc.) but carries its own date,
# volume and origin, matching the shape of every other row in the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: insert a blank row at 378, shifting the
# existing rows 378-408 down to 379-409.
$ws.Rows.Item(378).Insert()

$ws.Cells.Item(378, 1).Value = 10
$ws.Cells.Item(378, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(378, 3).Value = "La Araucanía"
$ws.Cells.Item(378, 4).Value = 44783
$ws.Cells.Item(378, 5).Value = 9
$ws.Cells.Item(378, 6).Value = 100112037
$ws.Cells.Item(378, 7).Value = "Cebollín"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 40
$ws.Cells.Item(378, 11).Value = 7000
$ws.Cells.Item(378, 12).Value = 7000
$ws.Cells.Item(378, 13).Value = 7000
$ws.Cells.Item(378, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(378, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(378, 16).Value = 583
$ws.Cells.Item(378, 17).Value = 12
$ws.Cells.Item(378, 18).Value = "Hortaliza"
